# Auto-generated Excel COM-interop script
# Applies the "Add data for 2022-10-22" update to violent-crime-full-year.xlsx
# For each affected worksheet, updates the 2022 (column I) values, plus a handful
# of minor corrections in columns D and H that were revised alongside it.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 5916
$ws.Range("I3").Value = 6154
$ws.Range("D4").Value = 1935
$ws.Range("H4").Value = 1673
$ws.Range("I4").Value = 1416
$ws.Range("I5").Value = 568
$ws.Range("I6").Value = 6955
$ws.Range("D7").Value = 28125
$ws.Range("H7").Value = 25984
$ws.Range("I7").Value = 21009

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 163
$ws.Range("I7").Value = 663
$ws.Range("I8").Value = 1261
$ws.Range("I11").Value = 314
$ws.Range("I14").Value = 118
$ws.Range("I15").Value = 235
$ws.Range("I18").Value = 157
$ws.Range("I19").Value = 582
$ws.Range("I20").Value = 524
$ws.Range("I23").Value = 211
$ws.Range("I25").Value = 111
$ws.Range("I26").Value = 28
$ws.Range("I29").Value = 1307
$ws.Range("I31").Value = 209
$ws.Range("I33").Value = 952
$ws.Range("I36").Value = 280
$ws.Range("I37").Value = 671
$ws.Range("I40").Value = 38
$ws.Range("I41").Value = 88
$ws.Range("I42").Value = 715
$ws.Range("I43").Value = 180
$ws.Range("I49").Value = 142
$ws.Range("I52").Value = 452
$ws.Range("I53").Value = 219
$ws.Range("I54").Value = 432
$ws.Range("I58").Value = 13
$ws.Range("I59").Value = 35
$ws.Range("D63").Value = 324
$ws.Range("I63").Value = 78
$ws.Range("I65").Value = 490
$ws.Range("I67").Value = 812
$ws.Range("I68").Value = 76
$ws.Range("I72").Value = 84
$ws.Range("I74").Value = 32
$ws.Range("I77").Value = 135
$ws.Range("I79").Value = 596
$ws.Range("H83").Value = 545
$ws.Range("I83").Value = 454
$ws.Range("I84").Value = 183
$ws.Range("I85").Value = 955
$ws.Range("I88").Value = 191
$ws.Range("I90").Value = 255
$ws.Range("I91").Value = 228
$ws.Range("I93").Value = 118
$ws.Range("I94").Value = 220
$ws.Range("I96").Value = 227
$ws.Range("D101").Value = 28125
$ws.Range("H101").Value = 25984
$ws.Range("I101").Value = 21009

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 271
$ws.Range("I3").Value = 367
$ws.Range("I5").Value = 33
$ws.Range("I7").Value = 955

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I2").Value = 124
$ws.Range("I4").Value = 37
$ws.Range("I6").Value = 117
$ws.Range("I7").Value = 452

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I2").Value = 131
$ws.Range("I6").Value = 83
$ws.Range("I7").Value = 314

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I3").Value = 356
$ws.Range("I4").Value = 75
$ws.Range("I5").Value = 36
$ws.Range("I6").Value = 408
$ws.Range("I7").Value = 1261

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I2").Value = 50
$ws.Range("I6").Value = 100
$ws.Range("I7").Value = 219

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I3").Value = 206
$ws.Range("I4").Value = 35
$ws.Range("I5").Value = 31
$ws.Range("I6").Value = 175
$ws.Range("I7").Value = 663

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I2").Value = 67
$ws.Range("I7").Value = 227

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("I6").Value = 42
$ws.Range("I7").Value = 118

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 202
$ws.Range("I3").Value = 223
$ws.Range("I7").Value = 671

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I3").Value = 299
$ws.Range("I5").Value = 24
$ws.Range("I7").Value = 812

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I6").Value = 83
$ws.Range("I7").Value = 209

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("I2").Value = 67
$ws.Range("I7").Value = 183

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I3").Value = 149
$ws.Range("I7").Value = 490

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I3").Value = 163
$ws.Range("H4").Value = 38
$ws.Range("I6").Value = 98
$ws.Range("H7").Value = 545
$ws.Range("I7").Value = 454

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I3").Value = 362
$ws.Range("I6").Value = 300
$ws.Range("I7").Value = 952

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("I4").Value = 16
$ws.Range("I7").Value = 142

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I6").Value = 206
$ws.Range("I7").Value = 432

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 382
$ws.Range("I3").Value = 451
$ws.Range("I6").Value = 361
$ws.Range("I7").Value = 1307

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I3").Value = 178
$ws.Range("I6").Value = 172
$ws.Range("I7").Value = 582

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("I6").Value = 20
$ws.Range("I7").Value = 88

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I2").Value = 181
$ws.Range("I4").Value = 50
$ws.Range("I6").Value = 226
$ws.Range("I7").Value = 715

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I6").Value = 61
$ws.Range("I7").Value = 211

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I6").Value = 62
$ws.Range("I7").Value = 228

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I2").Value = 173
$ws.Range("I3").Value = 192
$ws.Range("I6").Value = 174
$ws.Range("I7").Value = 596

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I3").Value = 150
$ws.Range("I6").Value = 182
$ws.Range("I7").Value = 524

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("I2").Value = 46
$ws.Range("I3").Value = 36
$ws.Range("I6").Value = 68
$ws.Range("I7").Value = 157

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I3").Value = 91
$ws.Range("I7").Value = 280

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("I4").Value = 6
$ws.Range("I7").Value = 118

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("I4").Value = 14
$ws.Range("I7").Value = 220

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("I2").Value = 42
$ws.Range("I7").Value = 111

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("I2").Value = 73
$ws.Range("I7").Value = 235

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("I6").Value = 18
$ws.Range("I7").Value = 28

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("I2").Value = 17
$ws.Range("I7").Value = 35

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I6").Value = 34
$ws.Range("I7").Value = 163

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("I3").Value = 68
$ws.Range("I7").Value = 191

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I6").Value = 85
$ws.Range("I7").Value = 255

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("I6").Value = 18
$ws.Range("I7").Value = 76

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I6").Value = 101
$ws.Range("I7").Value = 180

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("I2").Value = 18
$ws.Range("I7").Value = 84

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("I3").Value = 47
$ws.Range("I4").Value = 6
$ws.Range("I7").Value = 135

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("I3").Value = 16
$ws.Range("I7").Value = 38

$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("I6").Value = 22
$ws.Range("I7").Value = 32

$ws = $wb.Worksheets.Item("Millenium Park")
$ws.Range("I6").Value = 6
$ws.Range("I7").Value = 13

